$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SVGA (column E) timing values take over the old "WORKS" (column F) numbers ---
$ws.Range("E11").Value = 128
$ws.Range("E12").Value = 350
$ws.Range("E13").Value = 194
$ws.Range("E14").Value = 2592

# --- Column F ("WORKS") is retired: clear the label and the now-redundant numbers ---
$ws.Range("F9").ClearContents()
$ws.Range("F9").HorizontalAlignment = -4152   # xlRight
$ws.Range("F10:F14").ClearContents()

# --- Derived ratio rows (15-17) no longer carry formulas, just keep their formatting ---
$ws.Range("B15:E17").ClearContents()

# --- Column widths: B:D and E shrink slightly, new column F gets the same width ---
$ws.Columns("B:D").ColumnWidth = 11.3
$ws.Range("E1").ColumnWidth = 11.3
$ws.Range("F1").ColumnWidth = 11.3

# --- Selection moves to the ratio block ---
$ws.Range("B15:E17").Select()
